# Increase MaxInvest Storage Adapt Szenarios Existing Units
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F = "MaxInvest" -> raise the invest caps for the existing-unit rows
$ws.Range("F8").Value = 16
$ws.Range("F10").Value = 15
$ws.Range("F14").Value = 9
$ws.Range("F16").Value = 77

# Leave the cursor/selection where the author ended up (I9:I18)
$null = $ws.Range("I9:I18").Select()
